$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text (Volume/Number and date-range banner) ---
$ws.Range("A8").Value = "Volume 31   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/24/2024  Through  6/30/2024"

# --- Crime Complaints weekly table (rows 15-33) ---

# Row 15
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 0

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 66.666666666666
$ws.Range("L16").Value = 30.434782608695
$ws.Range("M16").Value = -9.090909090909
$ws.Range("N16").Value = -78.260869565217

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 91
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 75
$ws.Range("L17").Value = 54.237288135593
$ws.Range("M17").Value = 97.826086956521
$ws.Range("N17").Value = -9

# Row 18
$ws.Range("F18").Value = "'0"
$ws.Range("A18").Copy() | Out-Null
$ws.Range("F18").PasteSpecial(-4122)
$ws.Range("H18").Value = -100
$ws.Range("I18").Value = 23
$ws.Range("K18").Value = 53.333333333333
$ws.Range("L18").Value = -11.538461538461
$ws.Range("M18").Value = -20.689655172413
$ws.Range("N18").Value = -90.254237288135

# Row 19
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 9
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = -40
$ws.Range("I19").Value = 68
$ws.Range("J19").Value = 83
$ws.Range("K19").Value = -18.072289156626
$ws.Range("L19").Value = 1.492537313432
$ws.Range("M19").Value = 23.636363636363
$ws.Range("N19").Value = -19.047619047619

# Row 20
$ws.Range("C20").Value = "'0"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = "'0"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("A20").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 170
$ws.Range("N20").Value = -84.117647058823

# Row 21
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 9
$ws.Range("F21").Value = 35
$ws.Range("G21").Value = 42
$ws.Range("H21").Value = -16.666666666666
$ws.Range("I21").Value = 241
$ws.Range("J21").Value = 187
$ws.Range("K21").Value = 28.877005347593
$ws.Range("L21").Value = 15.311004784689
$ws.Range("M21").Value = 36.158192090395
$ws.Range("N21").Value = -67.344173441734

# Row 22
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E22").Value = -100
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 133.333333333333
$ws.Range("M22").Value = -30

# Row 23
$ws.Range("F23").Value = 3

# Row 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 142.857142857143
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = 40
$ws.Range("I24").Value = 212
$ws.Range("J24").Value = 220
$ws.Range("K24").Value = -3.636363636363
$ws.Range("L24").Value = -16.535433070866
$ws.Range("M24").Value = 60.606060606060

# Row 25
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 450
$ws.Range("F25").Value = 28
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 102
$ws.Range("J25").Value = 128
$ws.Range("K25").Value = -20.3125
$ws.Range("L25").Value = -15

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 15
$ws.Range("I26").Value = 127
$ws.Range("J26").Value = 108
$ws.Range("K26").Value = 17.592592592592
$ws.Range("L26").Value = 28.282828282828
$ws.Range("M26").Value = -30.219780219780

# Row 27
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 0

# Row 28
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("F28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 12
$ws.Range("K28").Value = -25
$ws.Range("L28").Value = -35.714285714285

# Row 29
$ws.Range("D29").Value = 1
$ws.Range("F29").Value = "'0"
$ws.Range("A29").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -22.222222222222

# Row 30
$ws.Range("F30").Value = "'0"
$ws.Range("A30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 0

# Row 33
$ws.Range("D33").NumberFormat = '#,##0'
$ws.Range("D33").Value = 1
$ws.Range("E33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E33").Value = -100
$ws.Range("G33").NumberFormat = '#,##0'
$ws.Range("G33").Value = 1
$ws.Range("H33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = -50
